# Refresh cryptos list values (coinranking.com scrape, Apr 24 2023 run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.662.28'

# Row 3
$ws.Range('D3').Value = '1.865.70'
$ws.Range('E3').Value = '  -0.81%  '

# Row 4
$ws.Range('E4').Value = '  +0.42%  '

# Row 5
$ws.Range('D5').Value = "'333.57"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.12%  '

# Row 6
$ws.Range('E6').Value = '  +0.23%  '

# Row 7
$ws.Range('D7').Value = "'0.4706"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.06%  '

# Row 8
$ws.Range('D8').Value = "'0.3927"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.24%  '

# Row 9
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = "'45.50"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.34%  '

# Row 10
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.07997"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.77%  '

# Row 11
$ws.Range('D11').Value = "'1.003"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.36%  '

# Row 12
$ws.Range('D12').Value = "'21.85"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.52%  '

# Row 13
$ws.Range('D13').Value = "'6.002"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.45%  '

# Row 14
$ws.Range('D14').Value = '1.861.99'
$ws.Range('E14').Value = '  -1.08%  '

# Row 15
$ws.Range('D15').Value = "'7.256"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.51%  '

# Row 16
$ws.Range('D16').Value = "'1.011"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.17%  '

# Row 17
$ws.Range('D17').Value = "'88.52"
$ws.Range('D17').Style = "Normal"

# Row 18
$ws.Range('D18').Value = "'0.06725"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.13%  '

# Row 19
$ws.Range('D19').Value = "'0.00001043"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.22%  '

# Row 20
$ws.Range('D20').Value = "'17.15"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.12%  '

# Row 21
$ws.Range('E21').Value = '  +0.32%  '

# Row 22
$ws.Range('D22').Value = '27.637.32'
$ws.Range('E22').Value = '  -0.37%  '

# Row 23
$ws.Range('D23').Value = "'5.464"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.39%  '

# Row 24
$ws.Range('D24').Value = "'10.92"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.98%  '

# Row 25
$ws.Range('D25').Value = "'2.315"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.08%  '

# Row 26
$ws.Range('D26').Value = '2.085.97'
$ws.Range('E26').Value = '  -1.03%  '

# Row 27
$ws.Range('D27').Value = "'159.21"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.59%  '

# Row 28
$ws.Range('D28').Value = "'19.81"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.94%  '

# Row 29
$ws.Range('D29').Value = "'2.157"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.31%  '

# Row 30
$ws.Range('D30').Value = "'5.453"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.51%  '

# Row 31
$ws.Range('D31').Value = "'121.86"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.04%  '

# Row 32
$ws.Range('D32').Value = "'0.9831"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.28%  '

# Row 33
$ws.Range('D33').Value = "'0.09498"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.20%  '

# Row 34
$ws.Range('E34').Value = '  -0.14%  '

# Row 35
$ws.Range('D35').Value = "'5.316"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.91%  '

# Row 36
$ws.Range('D36').Value = "'1.338"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -8.12%  '

# Row 37
$ws.Range('D37').Value = "'0.06061"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.31%  '

# Row 38
$ws.Range('D38').Value = "'0.02233"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.37%  '

# Row 39
$ws.Range('D39').Value = "'8.344"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.51%  '

# Row 40
$ws.Range('D40').Value = "'1.194"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.98%  '

# Row 41
$ws.Range('D41').Value = "'1.008"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.27%  '

# Row 42
$ws.Range('D42').Value = "'0.5983"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.38%  '

# Row 43
$ws.Range('D43').Value = "'0.1887"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.70%  '

# Row 44
$ws.Range('D44').Value = "'10.30"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.15%  '

# Row 45
$ws.Range('D45').Value = "'1.248"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.66%  '

# Row 46
$ws.Range('D46').Value = "'0.5657"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.28%  '

# Row 47
$ws.Range('D47').Value = "'12.21"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.17%  '

# Row 48
$ws.Range('D48').Value = "'1.924"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.13%  '

# Row 49
$ws.Range('D49').Value = "'0.06759"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.03%  '

# Row 50
$ws.Range('D50').Value = "'112.01"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.05%  '

# Row 51
$ws.Range('D51').Value = "'3.060"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -9.86%  '

